$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "work" item in column F next to row 5 (Functionality section),
# entered as a literal text string "=" so Excel stores it as text, not a formula.
$ws.Range("F5").Value = '"="'

# Center the new cell's text horizontally, matching the new cellXf (s="2")
# with <alignment horizontal="center"/>.
$ws.Range("F5").HorizontalAlignment = -4108

# Move the active selection to G5, as reflected in the updated <selection>.
$ws.Range("G5").Select()
